# Generate Report for Handoff
# Update the localization-status report with the latest handoff run results:
#  - Overview sheet: refresh "Latest HO Xliff Generate Date" timestamps
#  - zh-cn / de-de sheets: refresh "Latest Handoff Datetime" timestamps and
#    set the "Priority" column to "ht" for the rows that were just handed off.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 13)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-30 18:27:55"
}

# --- zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-30 18:27:50"
}

# --- de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-30 18:27:55"
}
